$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.663.87"
$ws.Range("E2").Value = "'  +0.63%  "
$ws.Range("D3").Value = "'1.953.36"
$ws.Range("E3").Value = "'  +2.13%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'247.74"
$ws.Range("E5").Value = "'  +1.18%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D7").Value = "'0.4806"
$ws.Range("E7").Value = "'  -0.24%  "
$ws.Range("B8").Value = "'OKB"
$ws.Range("C8").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'44.37"
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("B9").Value = "'Cardano"
$ws.Range("C9").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2937"
$ws.Range("E9").Value = "'  +1.70%  "
$ws.Range("B10").Value = "'Dogecoin"
$ws.Range("C10").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06803"
$ws.Range("E10").Value = "'  +1.24%  "
$ws.Range("B11").Value = "'Litecoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "'112.06"
$ws.Range("E11").Value = "'  +1.44%  "
$ws.Range("B12").Value = "'Solana"
$ws.Range("C12").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'19.42"
$ws.Range("E12").Value = "'  +2.15%  "
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.963.89"
$ws.Range("E13").Value = "'  +2.50%  "
$ws.Range("B14").Value = "'TRON"
$ws.Range("C14").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.07683"
$ws.Range("E14").Value = "'  +1.86%  "
$ws.Range("B15").Value = "'Polkadot"
$ws.Range("C15").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.482"
$ws.Range("E15").Value = "'  +4.08%  "
$ws.Range("B16").Value = "'Polygon"
$ws.Range("C16").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.6857"
$ws.Range("E16").Value = "'  +2.10%  "
$ws.Range("B17").Value = "'BitcoinCash"
$ws.Range("C17").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "'293.98"
$ws.Range("E17").Value = "'  +2.16%  "
$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'30.666.65"
$ws.Range("E18").Value = "'  +0.60%  "
$ws.Range("B19").Value = "'Avalanche"
$ws.Range("C19").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'13.22"
$ws.Range("E19").Value = "'  +2.93%  "
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'5.651"
$ws.Range("E20").Value = "'  +3.14%  "
$ws.Range("D21").Value = "'2.215.88"
$ws.Range("E21").Value = "'  +2.25%  "
$ws.Range("B22").Value = "'ShibaInu"
$ws.Range("C22").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "'0.000007675"
$ws.Range("E22").Value = "'  +1.35%  "
$ws.Range("B23").Value = "'Dai"
$ws.Range("C23").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "'  +0.07%  "
$ws.Range("B24").Value = "'BinanceUSD"
$ws.Range("C24").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'0.9993"
$ws.Range("E24").Value = "'  +0.13%  "
$ws.Range("B25").Value = "'Chainlink"
$ws.Range("C25").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'6.611"
$ws.Range("E25").Value = "'  +3.08%  "
$ws.Range("B26").Value = "'Cosmos"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.759"
$ws.Range("E26").Value = "'  +3.33%  "
$ws.Range("B27").Value = "'Monero"
$ws.Range("C27").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'168.83"
$ws.Range("E27").Value = "'  +2.67%  "
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.26"
$ws.Range("E28").Value = "'  -0.23%  "
$ws.Range("B29").Value = "'LidoDAOToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.189"
$ws.Range("E29").Value = "'  +3.38%  "
$ws.Range("B30").Value = "'Stellar"
$ws.Range("C30").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1088"
$ws.Range("E30").Value = "'  +3.30%  "
$ws.Range("B31").Value = "'Toncoin"
$ws.Range("C31").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.435"
$ws.Range("E31").Value = "'  +2.21%  "
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.676"
$ws.Range("E32").Value = "'  +15.86%  "
$ws.Range("B33").Value = "'InternetComputer(DFINITY)"
$ws.Range("C33").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.456"
$ws.Range("E33").Value = "'  +6.80%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.05064"
$ws.Range("E34").Value = "'  +1.75%  "
$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7748"
$ws.Range("E35").Value = "'  +6.16%  "
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.163"
$ws.Range("E36").Value = "'  +2.64%  "
$ws.Range("B37").Value = "'VeChain"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02075"
$ws.Range("E37").Value = "'  +2.23%  "
$ws.Range("B38").Value = "'HuobiToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.730"
$ws.Range("E38").Value = "'  +0.35%  "
$ws.Range("B39").Value = "'MXToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.697"
$ws.Range("E39").Value = "'  +1.11%  "
$ws.Range("B40").Value = "'RenderToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.059"
$ws.Range("E40").Value = "'  +1.97%  "
$ws.Range("B41").Value = "'Quant"
$ws.Range("C41").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'110.79"
$ws.Range("E41").Value = "'  +0.20%  "
$ws.Range("B42").Value = "'TheSandbox"
$ws.Range("C42").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.4463"
$ws.Range("E42").Value = "'  +0.45%  "
$ws.Range("B43").Value = "'TrustWalletToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8716"
$ws.Range("E43").Value = "'  +0.73%  "
$ws.Range("B44").Value = "'FraxShare"
$ws.Range("C44").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.983"
$ws.Range("E44").Value = "'  +3.41%  "
$ws.Range("B45").Value = "'PaxDollar"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "'  +0.03%  "
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'69.50"
$ws.Range("E46").Value = "'  +2.09%  "
$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.394"
$ws.Range("E47").Value = "'  +0.93%  "
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.312"
$ws.Range("B49").Value = "'Algorand"
$ws.Range("C49").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1253"
$ws.Range("E49").Value = "'  +1.09%  "
$ws.Range("B50").Value = "'BitcoinSV"
$ws.Range("C50").Value = "'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'47.90"
$ws.Range("E50").Value = "'  -2.21%  "
$ws.Range("B51").Value = "'Elrond"
$ws.Range("C51").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.69"
$ws.Range("E51").Value = "'  +2.56%  "
